$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text updates (coin names / links) - safe to assign directly
$textUpdates = @{
    'B6' = 'FTXToken'
    'C6' = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
    'B7' = 'MXToken'
    'C7' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'B8' = 'BTSEToken'
    'C8' = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
    'B9' = 'LiechtensteinCryptoassetsExchange'
    'C9' = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
    'B10' = 'WazirX'
    'C10' = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    'B11' = 'MandalaExchangeToken'
    'C11' = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    'B12' = 'BitrueCoin'
    'C12' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    'B13' = 'BitMartToken'
    'C13' = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    'B14' = 'BitForexToken'
    'C14' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    'B15' = 'TigerCash'
    'C15' = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
    'B16' = 'LEO'
    'C16' = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    'B17' = 'GateToken'
    'C17' = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
}

foreach ($ref in $textUpdates.Keys) {
    $ws.Range($ref).Value = $textUpdates[$ref]
}

# Numeric-looking text updates (prices / percentages) - must stay as Text cells,
# matching the original inlineStr cell type (no numeric coercion, no % conversion).
# Use a quote-prefix to force text entry, then reset the style back to Normal so
# no stray formatting (e.g. quotePrefix / @ number format) is left on the cell.
$numericTextUpdates = @{
    'D2' = '297.52'
    'E2' = '2.54%'
    'D3' = '41.51'
    'E3' = '2.90%'
    'D4' = '5.027'
    'E4' = '-0.41%'
    'D5' = '0.07544'
    'E5' = '3.24%'
    'D6' = '1.592'
    'E6' = '2.25%'
    'D7' = '0.9289'
    'E7' = '1.27%'
    'D8' = '2.423'
    'E8' = '1.09%'
    'D9' = '0.1200'
    'E9' = '2.77%'
    'D10' = '0.1844'
    'E10' = '6.60%'
    'D11' = '0.08919'
    'E11' = '2.55%'
    'D12' = '0.04047'
    'E12' = '-2.75%'
    'D13' = '0.1055'
    'E13' = '0.06%'
    'D14' = '0.001284'
    'E14' = '1.20%'
    'D15' = '0.005864'
    'E15' = '0.43%'
    'D16' = '3.340'
    'E16' = '-1.57%'
    'D17' = '4.374'
    'E17' = '2.29%'
    'D18' = '0.3310'
    'E18' = '1.04%'
    'D19' = '7.972'
    'E19' = '1.98%'
    'D20' = '0.1419'
    'E20' = '5.03%'
    'D22' = '0.04051'
    'E22' = '4.68%'
    'E23' = '-0.19%'
    'D24' = '0.004229'
    'E24' = '8.94%'
    'E25' = '-3.94%'
    'E26' = '0.01%'
    'E38' = '4.06%'
    'D39' = '0.05213'
    'E39' = '4.78%'
    'D40' = '0.006508'
    'E40' = '-3.17%'
    'D41' = '0.007818'
    'D42' = '0.1332'
    'E42' = '4.55%'
    'D43' = '0.007602'
    'E43' = '3.25%'
    'D44' = '0.007829'
    'E44' = '11.00%'
    'D45' = '0.3223'
    'E45' = '11.39%'
    'D46' = '0.00006668'
    'E46' = '3.60%'
    'D49' = '0.04627'
    'E49' = '55.17%'
}

foreach ($ref in $numericTextUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.Value = "'" + $numericTextUpdates[$ref]
    $cell.Style = "Normal"
}

